# Insert a new data row at row 146 (shifting existing rows 146-163 down to 147-164)
# and populate the new row with the data for the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 146; this shifts rows 146..163 down
# to 147..164 and copies formatting (including the date style) from the row above.
$ws.Rows.Item(146).Insert()

# Fill in the newly inserted row 146 with the new record's values.
$ws.Range("A146").Value = 1
$ws.Range("B146").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C146").Value = "Arica y Parinacota"
$ws.Range("D146").Value = 45142
$ws.Range("E146").Value = 15
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100102
$ws.Range("H146").Value = "Cítricos"
$ws.Range("I146").Value = 100102005
$ws.Range("J146").Value = "Naranja"
$ws.Range("K146").Value = "Navel Late"
$ws.Range("L146").Value = "Segunda"
$ws.Range("M146").Value = 350
$ws.Range("N146").Value = 800
$ws.Range("O146").Value = 858
$ws.Range("P146").Value = 825
$ws.Range("Q146").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R146").Value = "Región de Coquimbo"
$ws.Range("S146").Value = 825
$ws.Range("T146").Value = 1
